$d = $word.ActiveDocument

# Update the date line in the first paragraph
$d.Content.Find.Execute("2024-10-07 Monday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-10-08 Tuesday", 2)

$t = $d.Tables.Item(1)

function Replace-CellText($table, $row, $col, $oldText, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $newText, 2)
}

# Row 1
Replace-CellText $t 1 1 "74×47=3478" "78×32=2496"
Replace-CellText $t 1 2 "34×22=748"  "79×67=5293"
Replace-CellText $t 1 3 "33×21=693"  "88×86=7568"
Replace-CellText $t 1 4 "25×11=275"  "45×75=3375"
Replace-CellText $t 1 5 "29×20=580"  "55×87=4785"

# Row 5
Replace-CellText $t 5 1 "52×49=2548" "34×47=1598"
Replace-CellText $t 5 2 "90×12=1080" "76×47=3572"
Replace-CellText $t 5 3 "13×33=429"  "85×30=2550"
Replace-CellText $t 5 4 "94×36=3384" "90×49=4410"
Replace-CellText $t 5 5 "52×43=2236" "30×81=2430"

# Row 10
Replace-CellText $t 10 1 "95×64=6080" "29×18=522"
Replace-CellText $t 10 2 "70×46=3220" "44×64=2816"
Replace-CellText $t 10 3 "79×92=7268" "53×32=1696"
Replace-CellText $t 10 4 "44×34=1496" "81×45=3645"
Replace-CellText $t 10 5 "17×93=1581" "73×44=3212"

# Row 15
Replace-CellText $t 15 1 "26×49=1274" "45×99=4455"
Replace-CellText $t 15 2 "14×82=1148" "79×72=5688"
Replace-CellText $t 15 3 "84×29=2436" "93×29=2697"
Replace-CellText $t 15 4 "69×62=4278" "33×90=2970"
Replace-CellText $t 15 5 "21×50=1050" "20×19=380"

# Row 20
Replace-CellText $t 20 1 "15×26=390"  "77×13=1001"
Replace-CellText $t 20 2 "55×45=2475" "29×63=1827"
Replace-CellText $t 20 3 "53×35=1855" "11×65=715"
Replace-CellText $t 20 4 "47×89=4183" "93×19=1767"
Replace-CellText $t 20 5 "77×14=1078" "25×11=275"

Write-Host "All replacements complete"
